$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# The new row (40) is a copy of row 39's values, except column A (TimeStamp)
# gets a new timestamp value corresponding to 2025-05-02 10:40:52.
$newRow = 40
$srcRow = 39

# Copy style/format from the source row's cells so the new row looks identical.
$srcRange = $ws.Range($ws.Cells.Item($srcRow, 1), $ws.Cells.Item($srcRow, 15))
$dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 15))
$srcRange.Copy() | Out-Null
$dstRange.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Set the values for the new row.
$ws.Cells.Item($newRow, 1).Value = 45779.4450462963
$ws.Cells.Item($newRow, 2).Value = 10
$ws.Cells.Item($newRow, 3).Value = 6
$ws.Cells.Item($newRow, 4).Value = 299
$ws.Cells.Item($newRow, 5).Value = 474
$ws.Cells.Item($newRow, 6).Value = 449
$ws.Cells.Item($newRow, 7).Value = 534
$ws.Cells.Item($newRow, 8).Value = 4036
$ws.Cells.Item($newRow, 9).Value = 534
$ws.Cells.Item($newRow, 10).Value = 2570
$ws.Cells.Item($newRow, 11).Value = 257
$ws.Cells.Item($newRow, 12).Value = 494
$ws.Cells.Item($newRow, 13).Value = 30
$ws.Cells.Item($newRow, 14).Value = 4314
$ws.Cells.Item($newRow, 15).Value = 5402

$wb.Save()
